$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 1040
$ws.Range("D4").Value = 1039

$ws.Range("D4").Select()
